$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (tab title)
$ws.Name = "Through 2022-07-04"

# Update the label for the July row
$ws.Range("A8").Value = "July (through 07-04)"

# Update July row (row 8) values
$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 7
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 8
$ws.Range("G8").Value = 10
$ws.Range("H8").Value = 17
$ws.Range("I8").Value = 27

# Update Total row (row 9) values
$ws.Range("B9").Value = 131
$ws.Range("C9").Value = 254
$ws.Range("D9").Value = 397
$ws.Range("E9").Value = 362
$ws.Range("F9").Value = 259
$ws.Range("G9").Value = 482
$ws.Range("H9").Value = 777
$ws.Range("I9").Value = 833
